$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 913.7273
$ws.Range("I18").Value = 960.7
$ws.Range("J18").Value = 444
$ws.Range("K18").Value = 960.7
$ws.Range("L18").Value = 444
$ws.Range("M18").Value = -676.7
$ws.Range("N18").Value = -1012

$ws.Range("H55").Value = 59.4
$ws.Range("I55").Value = 72.26667
$ws.Range("J55").Value = 40.1
$ws.Range("K55").Value = 72.26667
$ws.Range("L55").Value = 40.1
$ws.Range("M55").Value = 141.73333
$ws.Range("N55").Value = -468.1

$ws.Range("H70").Value = 2047.4584
$ws.Range("I70").Value = 1621.875
$ws.Range("J70").Value = 2898.625
$ws.Range("K70").Value = 4865.625
$ws.Range("L70").Value = 8695.875
$ws.Range("M70").Value = -4595.625
$ws.Range("N70").Value = -9235.875

$ws.Range("H73").Value = 2047.4584
$ws.Range("I73").Value = 1621.875
$ws.Range("J73").Value = 2898.625
$ws.Range("K73").Value = 4865.625
$ws.Range("L73").Value = 8695.875
$ws.Range("M73").Value = -3929.625
$ws.Range("N73").Value = -10567.875

$ws.Range("H137").Value = 4069.7778
$ws.Range("I137").Value = 3825.6428
$ws.Range("J137").Value = 4924.25
$ws.Range("K137").Value = 11476.9284
$ws.Range("L137").Value = 14772.75
$ws.Range("M137").Value = -8926.928400000001
$ws.Range("N137").Value = -19872.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2335.4546
$ws.Range("I32").Value = 1254.9791
$ws.Range("J32").Value = 9744.429
$ws.Range("K32").Value = 1254.9791
$ws.Range("L32").Value = 9744.429
$ws.Range("M32").Value = -967.9791
$ws.Range("N32").Value = -10318.429

$ws.Range("H45").Value = 2135.2144
$ws.Range("I45").Value = 1900.6364
$ws.Range("K45").Value = 1900.6364
$ws.Range("M45").Value = -1523.6364

$ws.Range("H74").Value = 2449.2183
$ws.Range("I74").Value = 1541.0741
$ws.Range("J74").Value = 3324.9285
$ws.Range("K74").Value = 1541.0741
$ws.Range("L74").Value = 3324.9285
$ws.Range("M74").Value = -667.0741
$ws.Range("N74").Value = -5072.9285

$ws.Range("H77").Value = 2449.2183
$ws.Range("I77").Value = 1541.0741
$ws.Range("J77").Value = 3324.9285
$ws.Range("K77").Value = 7705.3705
$ws.Range("L77").Value = 16624.6425
$ws.Range("M77").Value = -3337.3705
$ws.Range("N77").Value = -25360.6425

$ws.Range("H102").Value = 30372.637
$ws.Range("I102").Value = 9497.6875
$ws.Range("K102").Value = 9497.6875
$ws.Range("M102").Value = -7875.6875

$ws.Range("H132").Value = 3667.6667
$ws.Range("I132").Value = 3304.6453
$ws.Range("J132").Value = 5918.4
$ws.Range("K132").Value = 9913.9359
$ws.Range("L132").Value = 17755.2
$ws.Range("M132").Value = -7383.9359
$ws.Range("N132").Value = -22815.2

$ws.Range("H139").Value = 55000
$ws.Range("J139").Value = 55000
$ws.Range("L139").Value = 55000
$ws.Range("N139").Value = -65280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 90910950
$ws.Range("I99").Value = 111112740
$ws.Range("J99").Value = 2949
$ws.Range("K99").Value = 111112740
$ws.Range("L99").Value = 2949
$ws.Range("M99").Value = -111111242
$ws.Range("N99").Value = -5945

$ws.Range("H134").Value = 3198.4546
$ws.Range("I134").Value = 2886.2
$ws.Range("K134").Value = 8658.599999999999
$ws.Range("M134").Value = -6123.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 39500
$ws.Range("J75").Value = 39500
$ws.Range("L75").Value = 39500
$ws.Range("N75").Value = -41496

$ws.Range("H78").Value = 39500
$ws.Range("J78").Value = 39500
$ws.Range("L78").Value = 118500
$ws.Range("N78").Value = -128484

$ws.Range("H134").Value = 3761.6667
$ws.Range("I134").Value = 3626.5557
$ws.Range("J134").Value = 4369.6665
$ws.Range("K134").Value = 10879.6671
$ws.Range("L134").Value = 13108.9995
$ws.Range("M134").Value = -8344.667099999999
$ws.Range("N134").Value = -18178.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 5000001.5
$ws.Range("I10").Value = 5000001.5
$ws.Range("K10").Value = 5000001.5
$ws.Range("M10").Value = -4999832.5

$ws.Range("H11").Value = 11682750
$ws.Range("I11").Value = 9108454
$ws.Range("K11").Value = 9108454
$ws.Range("M11").Value = -9108315

$ws.Range("H13").Value = 2000
$ws.Range("I13").Value = 2000
$ws.Range("K13").Value = 2000
$ws.Range("M13").Value = -1861

$ws.Range("H41").Value = 8448.25
$ws.Range("I41").Value = 8448.25
$ws.Range("K41").Value = 8448.25
$ws.Range("M41").Value = -8093.25

$ws.Range("H62").Value = 40085
$ws.Range("J62").Value = 40085
$ws.Range("L62").Value = 40085
$ws.Range("N62").Value = -41457

$ws.Range("H64").Value = 82900
$ws.Range("I64").Value = 82900
$ws.Range("K64").Value = 82900
$ws.Range("M64").Value = -82652

$ws.Range("H65").Value = 40085
$ws.Range("J65").Value = 40085
$ws.Range("L65").Value = 120255
$ws.Range("N65").Value = -127119

$ws.Range("H67").Value = 82900
$ws.Range("I67").Value = 82900
$ws.Range("K67").Value = 82900
$ws.Range("M67").Value = -82042

$ws.Range("H70").Value = 42409.43
$ws.Range("I70").Value = 78678.57000000001
$ws.Range("K70").Value = 78678.57000000001
$ws.Range("M70").Value = -78408.57000000001

$ws.Range("H73").Value = 42409.43
$ws.Range("I73").Value = 78678.57000000001
$ws.Range("K73").Value = 78678.57000000001
$ws.Range("M73").Value = -77742.57000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5916.3335
$ws.Range("I22").Value = 501
$ws.Range("J22").Value = 6999.4
$ws.Range("K22").Value = 501
$ws.Range("L22").Value = 6999.4
$ws.Range("M22").Value = -206
$ws.Range("N22").Value = -7589.4

$ws.Range("H27").Value = 5916.3335
$ws.Range("I27").Value = 501
$ws.Range("J27").Value = 6999.4
$ws.Range("K27").Value = 501
$ws.Range("L27").Value = 6999.4
$ws.Range("M27").Value = -394
$ws.Range("N27").Value = -7213.4

$ws.Range("H136").Value = 4413.4614
$ws.Range("I136").Value = 4571.25
$ws.Range("J136").Value = 4161
$ws.Range("K136").Value = 13713.75
$ws.Range("L136").Value = 12483
$ws.Range("M136").Value = -11163.75
$ws.Range("N136").Value = -17583

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4134.067
$ws.Range("I81").Value = 3317.6
$ws.Range("J81").Value = 5767
$ws.Range("K81").Value = 6635.2
$ws.Range("L81").Value = 11534
$ws.Range("M81").Value = -5574.2
$ws.Range("N81").Value = -13656

$ws.Range("H84").Value = 4134.067
$ws.Range("I84").Value = 3317.6
$ws.Range("J84").Value = 5767
$ws.Range("K84").Value = 33176
$ws.Range("L84").Value = 57670
$ws.Range("M84").Value = -27872
$ws.Range("N84").Value = -68278
